# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# "System"/"system" is listed first in the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$dnasrRows = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)

foreach ($r in $dnasrRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# Rows where "admin@admin.com, System" -> "System, admin@admin.com"
$adminRows = @(7, 33, 59)

foreach ($r in $adminRows) {
    $ws.Range("G$r").Value = "System, admin@admin.com"
}

# Rows where "system, System, backup@backdoor.com" -> "System, system, backup@backdoor.com"
$backdoorRows = @(2, 28, 54)

foreach ($r in $backdoorRows) {
    $ws.Range("G$r").Value = "System, system, backup@backdoor.com"
}

$wb.Save()
